$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (including the date number-format style) of the last
# existing data row (75) down onto the two new rows (76-77) so the new
# date cells pick up the same style index used by the rest of column A.
$ws.Range("A75:G75").Copy()
$ws.Range("A76:G77").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 76 - 2025-03-01
$ws.Range("A76").Value = 45717
$ws.Range("B76").Value = -0.416
$ws.Range("C76").Value = -0.256
$ws.Range("D76").Value = 0.213
$ws.Range("E76").Value = 0.251
$ws.Range("F76").Value = 0.065
$ws.Range("G76").Value = 79.21

# Row 77 - 2025-04-01
$ws.Range("A77").Value = 45748
$ws.Range("B77").Value = -0.063
$ws.Range("C77").Value = -0.307
$ws.Range("D77").Value = 0.217
$ws.Range("E77").Value = 0.126
$ws.Range("F77").Value = 0.185
$ws.Range("G77").Value = 79.56
